$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Two new weekly price rows are inserted right before the existing row 312,
# pushing the rest of the "Naranja" data block down by 2 rows
# (old A1:T326 dimension becomes A1:T328).
$ws.Range("A312:A313").EntireRow.Insert()

# New row 312: Naranja / Valencia / Primera, week of 44585
$ws.Cells.Item(312, 1).Value = 7
$ws.Cells.Item(312, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(312, 3).Value = "Ñuble"
$ws.Cells.Item(312, 4).Value = 44585
$ws.Cells.Item(312, 5).Value = 16
$ws.Cells.Item(312, 6).Value = "Fruta"
$ws.Cells.Item(312, 7).Value = 100102
$ws.Cells.Item(312, 8).Value = "Cítricos"
$ws.Cells.Item(312, 9).Value = 100102005
$ws.Cells.Item(312, 10).Value = "Naranja"
$ws.Cells.Item(312, 11).Value = "Valencia"
$ws.Cells.Item(312, 12).Value = "Primera"
$ws.Cells.Item(312, 13).Value = 120
$ws.Cells.Item(312, 14).Value = 9000
$ws.Cells.Item(312, 15).Value = 10000
$ws.Cells.Item(312, 16).Value = 9500
$ws.Cells.Item(312, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(312, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(312, 19).Value = 633
$ws.Cells.Item(312, 20).Value = 15

# New row 313: Naranja / Valencia / Segunda, week of 44585
$ws.Cells.Item(313, 1).Value = 7
$ws.Cells.Item(313, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(313, 3).Value = "Ñuble"
$ws.Cells.Item(313, 4).Value = 44585
$ws.Cells.Item(313, 5).Value = 16
$ws.Cells.Item(313, 6).Value = "Fruta"
$ws.Cells.Item(313, 7).Value = 100102
$ws.Cells.Item(313, 8).Value = "Cítricos"
$ws.Cells.Item(313, 9).Value = 100102005
$ws.Cells.Item(313, 10).Value = "Naranja"
$ws.Cells.Item(313, 11).Value = "Valencia"
$ws.Cells.Item(313, 12).Value = "Segunda"
$ws.Cells.Item(313, 13).Value = 50
$ws.Cells.Item(313, 14).Value = 8000
$ws.Cells.Item(313, 15).Value = 8000
$ws.Cells.Item(313, 16).Value = 8000
$ws.Cells.Item(313, 17).Value = "$/bandeja 15 kilos granel"
$ws.Cells.Item(313, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(313, 19).Value = 533
$ws.Cells.Item(313, 20).Value = 15
